# Sanity Semilla 6 - ajustes en clases de portabilidad prepago y postpago
#
# Row 5 (prepago "consulta_log" placeholders) becomes two distinct, real
# values: F5 -> "app", G5 -> "tvcan1020Sem_6". Write G5 first so the shared
# string for "tvcan1020Sem_6" is interned before "app"'s, matching the
# author's original save order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = "tvcan1020Sem_6"
$ws.Range("F5").Value = "app"

# The author's last click before saving landed on F7.
[void]$ws.Range("F7").Select()
